# Update the p-values table in Fig 3: replace the Cod/Hake trend-line
# p-values that no longer depend on country/species factors, leaving the
# cells that are unchanged (OHI economic/Hake = 0.78, Vulnerability/Cod =
# "<0.01") untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing end-of-cell marker so the text is replaced
    # in place rather than inserted before it.
    $textRange = $d.Range($r.Start, $r.End - 1)
    $textRange.Text = $newText
}

# GDP 2016
Set-CellValue $t 2 2 "0.62"   # Cod:  0.93 -> 0.62
Set-CellValue $t 2 3 "0.98"   # Hake: 0.45 -> 0.98

# OHI fisheries
Set-CellValue $t 3 2 "0.22"   # Cod:  0.27 -> 0.22
Set-CellValue $t 3 3 "0.59"   # Hake: 0.02 -> 0.59

# OHI economic
Set-CellValue $t 4 2 "0.96"   # Cod:  0.92 -> 0.96
# Hake (0.78) is unchanged

# Readiness
Set-CellValue $t 5 2 "0.12"   # Cod:  0.15 -> 0.12
Set-CellValue $t 5 3 "0.76"   # Hake: 0.14 -> 0.76

# Vulnerability
# Cod (<0.01) is unchanged
Set-CellValue $t 6 3 "0.16"   # Hake: 0.01 -> 0.16
